$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 99997.336
$ws.Range("J120").Value = 99997.336
$ws.Range("L120").Value = 99997.336
$ws.Range("N120").Value = -109673.336

$ws.Range("H137").Value = 9171.786
$ws.Range("I137").Value = 14894
$ws.Range("J137").Value = 3449.5715
$ws.Range("K137").Value = 44682
$ws.Range("L137").Value = 10348.7145
$ws.Range("M137").Value = -42132
$ws.Range("N137").Value = -15448.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 24162.846
$ws.Range("I28").Value = 4705.6665
$ws.Range("K28").Value = 4705.6665
$ws.Range("M28").Value = -4513.6665

$ws.Range("H53").Value = 11166.333
$ws.Range("I53").Value = 11166.333
$ws.Range("K53").Value = 11166.333
$ws.Range("M53").Value = -10484.333

$ws.Range("H61").Value = 6108.2354
$ws.Range("I61").Value = 13962.25
$ws.Range("J61").Value = 3691.6155
$ws.Range("K61").Value = 13962.25
$ws.Range("L61").Value = 3691.6155
$ws.Range("M61").Value = -13750.25
$ws.Range("N61").Value = -4115.6155

$ws.Range("H63").Value = 2263.4375
$ws.Range("I63").Value = 1724.2106
$ws.Range("K63").Value = 1724.2106
$ws.Range("M63").Value = -1038.2106

$ws.Range("H66").Value = 2263.4375
$ws.Range("I66").Value = 1724.2106
$ws.Range("K66").Value = 8621.053
$ws.Range("M66").Value = -5189.053

$ws.Range("H74").Value = 3527.875
$ws.Range("I74").Value = 2077.25
$ws.Range("J74").Value = 4978.5
$ws.Range("K74").Value = 2077.25
$ws.Range("L74").Value = 4978.5
$ws.Range("M74").Value = -1203.25
$ws.Range("N74").Value = -6726.5

$ws.Range("H77").Value = 3527.875
$ws.Range("I77").Value = 2077.25
$ws.Range("J77").Value = 4978.5
$ws.Range("K77").Value = 10386.25
$ws.Range("L77").Value = 24892.5
$ws.Range("M77").Value = -6018.25
$ws.Range("N77").Value = -33628.5

$ws.Range("H99").Value = 24162.846
$ws.Range("I99").Value = 4705.6665
$ws.Range("K99").Value = 4705.6665
$ws.Range("M99").Value = -1710.6665

$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

$ws.Range("H132").Value = 3516.4856
$ws.Range("I132").Value = 3454.0967
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 10362.2901
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -7832.2901
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 6108.2354
$ws.Range("I136").Value = 13962.25
$ws.Range("J136").Value = 3691.6155
$ws.Range("K136").Value = 41886.75
$ws.Range("L136").Value = 11074.8465
$ws.Range("M136").Value = -39336.75
$ws.Range("N136").Value = -16174.8465

$ws.Range("H141").Value = 94400
$ws.Range("J141").Value = 88800
$ws.Range("L141").Value = 88800
$ws.Range("N141").Value = -99160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3174.25
$ws.Range("I107").Value = 948
$ws.Range("J107").Value = 5400.5
$ws.Range("K107").Value = 948
$ws.Range("L107").Value = 5400.5
$ws.Range("M107").Value = 972
$ws.Range("N107").Value = -9240.5

$ws.Range("H111").Value = 41157.4
$ws.Range("J111").Value = 41157.4
$ws.Range("L111").Value = 41157.4
$ws.Range("N111").Value = -49337.4

$ws.Range("H112").Value = 194996
$ws.Range("J112").Value = 194996
$ws.Range("L112").Value = 194996
$ws.Range("N112").Value = -197950

$ws.Range("H134").Value = 2546.611
$ws.Range("I134").Value = 2181.7144
$ws.Range("J134").Value = 3823.75
$ws.Range("K134").Value = 6545.1432
$ws.Range("L134").Value = 11471.25
$ws.Range("M134").Value = -4010.1432
$ws.Range("N134").Value = -16541.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9495
$ws.Range("J16").Value = 9495
$ws.Range("L16").Value = 9495
$ws.Range("N16").Value = -10069

$ws.Range("H22").Value = 133.33333
$ws.Range("I22").Value = 134.77777
$ws.Range("K22").Value = 134.77777
$ws.Range("M22").Value = 215.22223

$ws.Range("H31").Value = 2575.4468
$ws.Range("I31").Value = 1513.0667
$ws.Range("J31").Value = 4450.2354
$ws.Range("K31").Value = 1513.0667
$ws.Range("L31").Value = 4450.2354
$ws.Range("M31").Value = -1218.0667
$ws.Range("N31").Value = -5040.2354

$ws.Range("H34").Value = 2575.4468
$ws.Range("I34").Value = 1513.0667
$ws.Range("J34").Value = 4450.2354
$ws.Range("K34").Value = 1513.0667
$ws.Range("L34").Value = 4450.2354
$ws.Range("M34").Value = -1311.0667
$ws.Range("N34").Value = -4854.2354

$ws.Range("H113").Value = 9495
$ws.Range("J113").Value = 9495
$ws.Range("L113").Value = 9495
$ws.Range("N113").Value = -13835

$ws.Range("H122").Value = 15712.714
$ws.Range("I122").Value = 21317.8
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 63953.39999999999
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -61503.39999999999
$ws.Range("N122").Value = -10000

$ws.Range("H132").Value = 10280.704
$ws.Range("I132").Value = 3618.9285
$ws.Range("J132").Value = 17454.924
$ws.Range("K132").Value = 10856.7855
$ws.Range("L132").Value = 52364.772
$ws.Range("M132").Value = -8326.7855
$ws.Range("N132").Value = -57424.772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 511.2857
$ws.Range("I44").Value = 104.833336
$ws.Range("K44").Value = 314.500008
$ws.Range("M44").Value = 83.49999200000002

$ws.Range("H88").Value = 8567.286
$ws.Range("I88").Value = 6988.5
$ws.Range("J88").Value = 9198.799999999999
$ws.Range("K88").Value = 20965.5
$ws.Range("L88").Value = 27596.4
$ws.Range("M88").Value = -20537.5
$ws.Range("N88").Value = -28452.4

$ws.Range("H91").Value = 8567.286
$ws.Range("I91").Value = 6988.5
$ws.Range("J91").Value = 9198.799999999999
$ws.Range("K91").Value = 20965.5
$ws.Range("L91").Value = 27596.4
$ws.Range("M91").Value = -19483.5
$ws.Range("N91").Value = -30560.4

$ws.Range("H94").Value = 18175992
$ws.Range("I94").Value = 3500
$ws.Range("K94").Value = 10500
$ws.Range("M94").Value = -9824

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 130267
$ws.Range("J42").Value = 130267
$ws.Range("L42").Value = 130267
$ws.Range("N42").Value = -131237

$ws.Range("H70").Value = 6151.7
$ws.Range("I70").Value = 5443
$ws.Range("K70").Value = 5443
$ws.Range("M70").Value = -5173

$ws.Range("H73").Value = 6151.7
$ws.Range("I73").Value = 5443
$ws.Range("K73").Value = 5443
$ws.Range("M73").Value = -4507

$ws.Range("H97").Value = 540.0833
$ws.Range("I97").Value = 488.9
$ws.Range("J97").Value = 796
$ws.Range("K97").Value = 488.9
$ws.Range("L97").Value = 796
$ws.Range("M97").Value = 7.100000000000023
$ws.Range("N97").Value = -1788

$ws.Range("H107").Value = 948.7857
$ws.Range("I107").Value = 948.7857
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 948.7857
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 971.2143
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 3876.8
$ws.Range("I113").Value = 3346.25
$ws.Range("K113").Value = 3346.25
$ws.Range("M113").Value = -1176.25

$ws.Range("H115").Value = 130267
$ws.Range("J115").Value = 130267
$ws.Range("L115").Value = 130267
$ws.Range("N115").Value = -132617

$ws.Range("H127").Value = 70256
$ws.Range("J127").Value = 70256
$ws.Range("L127").Value = 70256
$ws.Range("N127").Value = -80176

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 4832.6665
$ws.Range("J105").Value = 4832.6665
$ws.Range("L105").Value = 4832.6665
$ws.Range("N105").Value = -11820.6665

$ws.Range("H121").Value = 156000
$ws.Range("J121").Value = 156000
$ws.Range("L121").Value = 156000
$ws.Range("N121").Value = -159494

$ws.Range("H125").Value = 96309.664
$ws.Range("J125").Value = 96309.664
$ws.Range("L125").Value = 96309.664
$ws.Range("N125").Value = -106149.664

$ws.Range("H131").Value = 64575
$ws.Range("J131").Value = 72333.336
$ws.Range("L131").Value = 72333.336
$ws.Range("N131").Value = -82413.336

$ws.Range("H136").Value = 7891.4346
$ws.Range("J136").Value = 4078.5715
$ws.Range("L136").Value = 12235.7145
$ws.Range("N136").Value = -17335.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 36358.75
$ws.Range("I70").Value = 35665.5
$ws.Range("K70").Value = 35665.5
$ws.Range("M70").Value = -35350.5

$ws.Range("H73").Value = 36358.75
$ws.Range("I73").Value = 35665.5
$ws.Range("K73").Value = 35665.5
$ws.Range("M73").Value = -34573.5

$ws.Range("H105").Value = 34083
$ws.Range("J105").Value = 34083
$ws.Range("L105").Value = 34083
$ws.Range("N105").Value = -41071

$ws.Range("H132").Value = 4265
$ws.Range("I132").Value = 4154.362
$ws.Range("K132").Value = 12463.086
$ws.Range("M132").Value = -9933.085999999999
